$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.680.27'
$ws.Range('E2').Value = '  +0.37%  '

# Row 3
$ws.Range('D3').Value = '3.504.86'
$ws.Range('E3').Value = '  -0.55%  '

# Row 4
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').Value = '''606.28'
$ws.Range('E5').Value = '  -1.24%  '

# Row 6
$ws.Range('D6').Value = '''152.08'
$ws.Range('E6').Value = '  +0.29%  '

# Row 7
$ws.Range('D7').Value = '3.502.01'
$ws.Range('E7').Value = '  -0.62%  '

# Row 8
$ws.Range('E8').Value = '  +0.03%  '

# Row 9
$ws.Range('D9').Value = '''0.489'
$ws.Range('E9').Value = '  +1.98%  '

# Row 10
$ws.Range('E10').Value = '  +2.78%  '

# Row 11
$ws.Range('D11').Value = '''7.62'
$ws.Range('E11').Value = '  +7.06%  '

# Row 12
$ws.Range('E12').Value = '  +1.34%  '

# Row 13
$ws.Range('E13').Value = '  -1.81%  '

# Row 14
$ws.Range('D14').Value = '''32.40'
$ws.Range('E14').Value = '  +0.77%  '

# Row 15
$ws.Range('D15').Value = '4.095.06'
$ws.Range('E15').Value = '  -0.47%  '

# Row 16
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '67.552.46'
$ws.Range('E16').Value = '  +0.24%  '

# Row 17
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.490.54'
$ws.Range('E17').Value = '  -1.11%  '

# Row 18
$ws.Range('E18').Value = '  -0.58%  '

# Row 19
$ws.Range('D19').Value = '''6.52'
$ws.Range('E19').Value = '  +1.95%  '

# Row 20
$ws.Range('E20').Value = '  +1.70%  '

# Row 21
$ws.Range('D21').Value = '''9.88'
$ws.Range('E21').Value = '  +4.28%  '

# Row 22
$ws.Range('D22').Value = '''447.14'
$ws.Range('E22').Value = '  +0.43%  '

# Row 23
$ws.Range('D23').Value = '''0.630'
$ws.Range('E23').Value = '  +0.67%  '

# Row 24
$ws.Range('D24').Value = '''78.31'
$ws.Range('E24').Value = '  +1.18%  '

# Row 25
$ws.Range('D25').Value = '3.643.25'
$ws.Range('E25').Value = '  -0.46%  '

# Row 26
$ws.Range('E26').Value = '  +0.01%  '

# Row 27
$ws.Range('D27').Value = '''0.0000126'
$ws.Range('E27').Value = '  -4.39%  '

# Row 28
$ws.Range('D28').Value = '''10.09'
$ws.Range('E28').Value = '  -1.82%  '

# Row 29
$ws.Range('D29').Value = '''8.76'
$ws.Range('E29').Value = '  +4.42%  '

# Row 30
$ws.Range('D30').Value = '''2.51'
$ws.Range('E30').Value = '  +0.36%  '

# Row 31
$ws.Range('D31').Value = '''1.64'
$ws.Range('E31').Value = '  +5.85%  '

# Row 32
$ws.Range('D32').Value = '''0.171'
$ws.Range('E32').Value = '  +3.82%  '

# Row 33
$ws.Range('D33').Value = '''1.00'
$ws.Range('E33').Value = '  +0.40%  '

# Row 34
$ws.Range('D34').Value = '''25.64'
$ws.Range('E34').Value = '  -1.07%  '

# Row 35
$ws.Range('D35').Value = '''6.16'
$ws.Range('E35').Value = '  -0.28%  '

# Row 36
$ws.Range('E36').Value = '  +0.82%  '

# Row 37
$ws.Range('D37').Value = '3.493.70'
$ws.Range('E37').Value = '  -0.45%  '

# Row 38
$ws.Range('E38').Value = '  -0.62%  '

# Row 39
$ws.Range('E39').Value = '  -0.03%  '

# Row 40
$ws.Range('D40').Value = '''2.32'
$ws.Range('E40').Value = '  +7.46%  '

# Row 41
$ws.Range('D41').Value = '''177.64'
$ws.Range('E41').Value = '  +0.00%  '

# Row 42
$ws.Range('D42').Value = '''0.999'
$ws.Range('E42').Value = '  +0.03%  '

# Row 43
$ws.Range('D43').Value = '''0.0894'
$ws.Range('E43').Value = '  +1.22%  '

# Row 44
$ws.Range('E44').Value = '  +0.36%  '

# Row 45
$ws.Range('D45').Value = '''0.893'
$ws.Range('E45').Value = '  +1.20%  '

# Row 46
$ws.Range('D46').Value = '''30.28'
$ws.Range('E46').Value = '  +6.57%  '

# Row 47
$ws.Range('D47').Value = '''46.45'
$ws.Range('E47').Value = '  +2.96%  '

# Row 48
$ws.Range('E48').Value = '  +3.18%  '

# Row 49
$ws.Range('E49').Value = '  -2.36%  '

# Row 50
$ws.Range('E50').Value = '  +0.35%  '

# Row 51
$ws.Range('E51').Value = '  +2.16%  '
